$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (shifts H.. onward one column to the right,
# Excel automatically updates all formulas / shared-formula refs / column widths).
$ws.Columns("H").Insert()

# Populate the newly inserted column H with a "zveranl" header and boolean (FALSE)
# values for the five data rows.
$ws.Range("H1").Value = "zveranl"
$ws.Range("H2").Value = $false
$ws.Range("H3").Value = $false
$ws.Range("H4").Value = $false
$ws.Range("H5").Value = $false
$ws.Range("H6").Value = $false

# Reset the view so it is no longer scrolled to show column H first, and move the
# active selection to H7.
$ws.Application.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("H7").Select()
